# Apply the "include ragwitz and van den berg" edit:
# Insert two new rows of historic-table data (source "Ragwitz et al. (2023)",
# DE, scenario "Nachfrage+Tech") into the boxplot historic table:
#   - one into the "Energy" / "final energy demand per capita and year" block
#   - one into the "Industry" / "final energy demand per capita and year | industry" block
# Every row below each insertion point shifts down by one, which is exactly
# what Excel's Rows.Insert() does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 11 (Energy block) and populate it ---
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value  = "Ragwitz et al. (2023)"
$ws.Cells.Item(11, 2).Value  = "DE"
$ws.Cells.Item(11, 3).Value  = "Nachfrage+Tech"
$ws.Cells.Item(11, 4).Value  = "Energy"
$ws.Cells.Item(11, 5).Value  = "final energy demand per capita and year"
$ws.Cells.Item(11, 6).Value  = 2045
$ws.Cells.Item(11, 7).Value  = "GJ/cap/year"
$ws.Cells.Item(11, 8).Value  = 60.08
$ws.Cells.Item(11, 9).Value  = 102.96
$ws.Cells.Item(11, 10).Value = 0.5835275835275835
$ws.Cells.Item(11, 11).Value = 439

# --- Insert new row 40 (Industry block) and populate it ---
$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value  = "Ragwitz et al. (2023)"
$ws.Cells.Item(40, 2).Value  = "DE"
$ws.Cells.Item(40, 3).Value  = "Nachfrage+Tech"
$ws.Cells.Item(40, 4).Value  = "Industry"
$ws.Cells.Item(40, 5).Value  = "final energy demand per capita and year | industry"
$ws.Cells.Item(40, 6).Value  = 2045
$ws.Cells.Item(40, 7).Value  = "GJ/cap/year"
$ws.Cells.Item(40, 8).Value  = 24.47
$ws.Cells.Item(40, 9).Value  = 29.12
$ws.Cells.Item(40, 10).Value = 0.840315934065934
$ws.Cells.Item(40, 11).Value = 438
